# Update cryptocurrency price/volume data (refresh from source feed).
# Rows 29/30 additionally swap which coin occupies each rank slot
# (Bittensor <-> Binance-PegBSC-USD), carrying their own price/link data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.204.59'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '2.662.32'
$ws.Range("E3").Value = '  +2.85%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'608.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.28%  '
$ws.Range("D6").Value = "'143.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  -1.03%  '
$ws.Range("D9").Value = '2.662.90'
$ws.Range("E9").Value = '  +2.84%  '
$ws.Range("E10").Value = '  -0.16%  '
$ws.Range("E11").Value = '  +1.01%  '
$ws.Range("E12").Value = '  +0.67%  '
$ws.Range("E13").Value = '  +1.66%  '
$ws.Range("D14").Value = "'27.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("D15").Value = '3.141.05'
$ws.Range("E15").Value = '  +3.00%  '
$ws.Range("D16").Value = '63.088.37'
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '2.662.30'
$ws.Range("E18").Value = '  +2.74%  '
$ws.Range("D19").Value = "'11.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.55%  '
$ws.Range("D20").Value = "'339.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.60%  '
$ws.Range("E21").Value = '  +0.89%  '
$ws.Range("E22").Value = '  +2.86%  '
$ws.Range("D24").Value = "'67.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("E25").Value = '  +2.23%  '
$ws.Range("E26").Value = '  -2.88%  '
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("D28").Value = "'8.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.19%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").Value = "'539.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +16.29%  '
$ws.Range("D31").Value = "'7.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.27%  '
$ws.Range("E32").Value = '  +5.30%  '
$ws.Range("E33").Value = '  +7.06%  '
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("D35").Value = "'172.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.52%  '
$ws.Range("E36").Value = '  +12.38%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D39").Value = "'19.21"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.48%  '
$ws.Range("D40").Value = "'1.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.48%  '
$ws.Range("D41").Value = "'174.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.24%  '
$ws.Range("E43").Value = '  +0.82%  '
$ws.Range("E44").Value = '  +3.01%  '
$ws.Range("E46").Value = '  -0.74%  '
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  +1.54%  '
$ws.Range("D49").Value = "'18.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.35%  '
$ws.Range("E50").Value = '  +2.45%  '
$ws.Range("D51").Value = "'11.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.56%  '
